$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Columns whose values rotate cyclically among rows 2, 4 and 5:
# new row2 = old row4, new row4 = old row5, new row5 = old row2
$cols = @("A","B","E","F","G","H","Q","R","S","Z","AB","AW","AX")

# Capture the "before" values for each of the three rows first,
# since rows 2, 4 and 5 all change and we must not read already-overwritten data.
$row2 = @{}
$row4 = @{}
$row5 = @{}
foreach ($col in $cols) {
    $row2[$col] = $ws.Range($col + "2").Value2
    $row4[$col] = $ws.Range($col + "4").Value2
    $row5[$col] = $ws.Range($col + "5").Value2
}

foreach ($col in $cols) {
    $ws.Range($col + "2").Value2 = $row4[$col]
    $ws.Range($col + "4").Value2 = $row5[$col]
    $ws.Range($col + "5").Value2 = $row2[$col]
}
